$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.746.19'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.091.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.12%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.078.52'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.20%  '

$ws.Range("E10").Value = '  -7.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.55'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.477'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000226'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.592.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.809.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.69%  '

$ws.Range("E17").Value = '  -1.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.091.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.13%  '

$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '504.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.92%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.700'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -8.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.28%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -13.75%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.06%  '

$ws.Range("E30").Value = '  -0.17%  '

$ws.Range("E31").Value = '  +1.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '518.25'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.56%  '

$ws.Range("E35").Value = '  -2.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.77%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0411'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.61%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.127'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.25%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0817'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.51'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.971.11'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.53%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.252'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.00%  '

$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0552'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.66%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.93%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.111'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.14%  '

$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.11%  '
